{"js": "const pairs = [\n  [\"16+63=\", \"29+59=\"],\n  [\"85-29=\", \"14+0=\"],\n  [\"14+75=\", \"47-19=\"],\n  [\"53-18=\", \"74-25=\"],\n  [\"63+6=\", \"94-2=\"],\n  [\"94-42=\", \"7+73=\"],\n  [\"80-8=\", \"70-50=\"],\n  [\"93-10=\", \"22+21=\"],\n  [\"51-24=\", \"10+89=\"],\n  [\"7+70=\", \"92-82=\"],\n  [\"78-45=\", \"71-49=\"],\n  [\"12+66=\", \"68-21=\"],\n  [\"61-48=\", \"15+62=\"],\n  [\"19+59=\", \"14+85=\"],\n  [\"91+8=\", \"82-80=\"],\n  [\"96+3=\", \"54+22=\"],\n  [\"67+28=\", \"9+53=\"],\n  [\"73-6=\", \"16+61=\"],\n  [\"18+37=\", \"23+34=\"],\n  [\"75+7=\", \"68+20=\"],\n  [\"38-11=\", \"41+2=\"],\n  [\"35-1=\", \"61-3=\"],\n  [\"41+52=\", \"2+2=\"],\n  [\"83-19=\", \"16+54=\"],\n  [\"94-66=\", \"38+25=\"],\n  [\"86-63=\", \"36+23=\"],\n  [\"86+4=\", \"26+50=\"],\n  [\"31-0=\", \"41+41=\"],\n  [\"79-32=\", \"25+19=\"],\n  [\"10+51=\", \"18+26=\"],\n  [\"49+34=\", \"50+10=\"],\n  [\"31-17=\", \"86-25=\"],\n  [\"60-29=\", \"80+15=\"],\n  [\"23+70=\", \"99-86=\"],\n  [\"54+16=\", \"44+34=\"],\n  [\"77-60=\", \"54+6=\"],\n  [\"84-41=\", \"77-5=\"],\n  [\"2+91=\", \"50+9=\"],\n  [\"41+6=\", \"94-77=\"],\n  [\"35+43=\", \"93-49=\"],\n  [\"90+4=\", \"48-7=\"],\n  [\"4+71=\", \"69-10=\"],\n  [\"79-9=\", \"67-60=\"],\n  [\"56-36=\", \"38+52=\"],\n  [\"62-0=\", \"73-55=\"],\n  [\"67+14=\", \"44+39=\"],\n  [\"18+5=\", \"2+39=\"],\n  [\"12+58=\", \"78+15=\"],\n  [\"35-19=\", \"2+78=\"],\n  [\"78+11=\", \"5+1=\"],\n  [\"47-4=\", \"92-55=\"],\n  [\"60-6=\", \"15+3=\"],\n  [\"8+34=\", \"14+14=\"],\n  [\"21+70=\", \"78+3=\"],\n  [\"42-13=\", \"27+34=\"],\n  [\"44-5=\", \"50-21=\"],\n  [\"28+13=\", \"84-63=\"],\n  [\"10+46=\", \"84-29=\"],\n  [\"85-59=\", \"10+65=\"],\n  [\"20+18=\", \"83-23=\"],\n  [\"69-6=\", \"61-24=\"],\n  [\"14+48=\", \"33+46=\"],\n  [\"77-61=\", \"68-43=\"],\n  [\"10+44=\", \"67-64=\"],\n  [\"53-2=\", \"45+35=\"],\n  [\"13+85=\", \"60-25=\"],\n  [\"78-24=\", \"23+28=\"],\n  [\"75-32=\", \"37-32=\"],\n  [\"23-6=\", \"54+10=\"],\n  [\"46-37=\", \"0+3=\"],\n  [\"27-2=\", \"48-13=\"],\n  [\"63-9=\", \"18-13=\"],\n  [\"32+22=\", \"88-26=\"],\n  [\"17+48=\", \"2+29=\"],\n  [\"50+29=\", \"82-60=\"],\n  [\"94-93=\", \"10-6=\"],\n  [\"68-14=\", \"83-76=\"],\n  [\"2+8=\", \"89-87=\"],\n  [\"92-31=\", \"35+52=\"],\n  [\"61-37=\", \"74-35=\"],\n  [\"42+7=\", \"79+13=\"],\n  [\"58+26=\", \"70+5=\"],\n  [\"26+5=\", \"25+62=\"],\n  [\"28+27=\", \"90-2=\"],\n  [\"33+54=\", \"51+30=\"],\n  [\"71-63=\", \"18+2=\"],\n  [\"75-26=\", \"3+13=\"],\n  [\"57+37=\", \"79-19=\"],\n  [\"50-17=\", \"77-6=\"],\n  [\"94-21=\", \"55+28=\"],\n  [\"21+10=\", \"83-1=\"],\n  [\"72-13=\", \"17+77=\"],\n  [\"15-2=\", \"54+13=\"],\n  [\"32+44=\", \"79+20=\"],\n  [\"40+58=\", \"46-27=\"],\n  [\"71-38=\", \"22-12=\"],\n  [\"7-7=\", \"16-16=\"],\n  [\"31+3=\", \"18+25=\"],\n  [\"5+27=\", \"18+65=\"],\n  [\"36+42=\", \"25+29=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"16+63=\", \"29+59=\"),\n    @(\"85-29=\", \"14+0=\"),\n    @(\"14+75=\", \"47-19=\"),\n    @(\"53-18=\", \"74-25=\"),\n    @(\"63+6=\", \"94-2=\"),\n    @(\"94-42=\", \"7+73=\"),\n    @(\"80-8=\", \"70-50=\"),\n    @(\"93-10=\", \"22+21=\"),\n    @(\"51-24=\", \"10+89=\"),\n    @(\"7+70=\", \"92-82=\"),\n    @(\"78-45=\", \"71-49=\"),\n    @(\"12+66=\", \"68-21=\"),\n    @(\"61-48=\", \"15+62=\"),\n    @(\"19+59=\", \"14+85=\"),\n    @(\"91+8=\", \"82-80=\"),\n    @(\"96+3=\", \"54+22=\"),\n    @(\"67+28=\", \"9+53=\"),\n    @(\"73-6=\", \"16+61=\"),\n    @(\"18+37=\", \"23+34=\"),\n    @(\"75+7=\", \"68+20=\"),\n    @(\"38-11=\", \"41+2=\"),\n    @(\"35-1=\", \"61-3=\"),\n    @(\"41+52=\", \"2+2=\"),\n    @(\"83-19=\", \"16+54=\"),\n    @(\"94-66=\", \"38+25=\"),\n    @(\"86-63=\", \"36+23=\"),\n    @(\"86+4=\", \"26+50=\"),\n    @(\"31-0=\", \"41+41=\"),\n    @(\"79-32=\", \"25+19=\"),\n    @(\"10+51=\", \"18+26=\"),\n    @(\"49+34=\", \"50+10=\"),\n    @(\"31-17=\", \"86-25=\"),\n    @(\"60-29=\", \"80+15=\"),\n    @(\"23+70=\", \"99-86=\"),\n    @(\"54+16=\", \"44+34=\"),\n    @(\"77-60=\", \"54+6=\"),\n    @(\"84-41=\", \"77-5=\"),\n    @(\"2+91=\", \"50+9=\"),\n    @(\"41+6=\", \"94-77=\"),\n    @(\"35+43=\", \"93-49=\"),\n    @(\"90+4=\", \"48-7=\"),\n    @(\"4+71=\", \"69-10=\"),\n    @(\"79-9=\", \"67-60=\"),\n    @(\"56-36=\", \"38+52=\"),\n    @(\"62-0=\", \"73-55=\"),\n    @(\"67+14=\", \"44+39=\"),\n    @(\"18+5=\", \"2+39=\"),\n    @(\"12+58=\", \"78+15=\"),\n    @(\"35-19=\", \"2+78=\"),\n    @(\"78+11=\", \"5+1=\"),\n    @(\"47-4=\", \"92-55=\"),\n    @(\"60-6=\", \"15+3=\"),\n    @(\"8+34=\", \"14+14=\"),\n    @(\"21+70=\", \"78+3=\"),\n    @(\"42-13=\", \"27+34=\"),\n    @(\"44-5=\", \"50-21=\"),\n    @(\"28+13=\", \"84-63=\"),\n    @(\"10+46=\", \"84-29=\"),\n    @(\"85-59=\", \"10+65=\"),\n    @(\"20+18=\", \"83-23=\"),\n    @(\"69-6=\", \"61-24=\"),\n    @(\"14+48=\", \"33+46=\"),\n    @(\"77-61=\", \"68-43=\"),\n    @(\"10+44=\", \"67-64=\"),\n    @(\"53-2=\", \"45+35=\"),\n    @(\"13+85=\", \"60-25=\"),\n    @(\"78-24=\", \"23+28=\"),\n    @(\"75-32=\", \"37-32=\"),\n    @(\"23-6=\", \"54+10=\"),\n    @(\"46-37=\", \"0+3=\"),\n    @(\"27-2=\", \"48-13=\"),\n    @(\"63-9=\", \"18-13=\"),\n    @(\"32+22=\", \"88-26=\"),\n    @(\"17+48=\", \"2+29=\"),\n    @(\"50+29=\", \"82-60=\"),\n    @(\"94-93=\", \"10-6=\"),\n    @(\"68-14=\", \"83-76=\"),\n    @(\"2+8=\", \"89-87=\"),\n    @(\"92-31=\", \"35+52=\"),\n    @(\"61-37=\", \"74-35=\"),\n    @(\"42+7=\", \"79+13=\"),\n    @(\"58+26=\", \"70+5=\"),\n    @(\"26+5=\", \"25+62=\"),\n    @(\"28+27=\", \"90-2=\"),\n    @(\"33+54=\", \"51+30=\"),\n    @(\"71-63=\", \"18+2=\"),\n    @(\"75-26=\", \"3+13=\"),\n    @(\"57+37=\", \"79-19=\"),\n    @(\"50-17=\", \"77-6=\"),\n    @(\"94-21=\", \"55+28=\"),\n    @(\"21+10=\", \"83-1=\"),\n    @(\"72-13=\", \"17+77=\"),\n    @(\"15-2=\", \"54+13=\"),\n    @(\"32+44=\", \"79+20=\"),\n    @(\"40+58=\", \"46-27=\"),\n    @(\"71-38=\", \"22-12=\"),\n    @(\"7-7=\", \"16-16=\"),\n    @(\"31+3=\", \"18+25=\"),\n    @(\"5+27=\", \"18+65=\"),\n    @(\"36+42=\", \"25+29=\"),\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
